# Apply "Natmi following Dr Hou advice" edit:
# Rebuild the LR-pair data rows as a full 3x3 Sending-cluster x Target-cluster
# grid (ECs, FAPs, sCs) for the fixed Adam23 -> Itga4 ligand-receptor pair,
# extending the sheet from 4 data rows to 9 data rows (rows 2-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Adam23-Itga4)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Adam23"
$ws.Cells.Item(2, 3).Value = "Itga4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2.0
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.03926266666666667
$ws.Cells.Item(2, 8).Value = 0.117788
$ws.Cells.Item(2, 9).Value = 0.005313231574131687
$ws.Cells.Item(2, 10).Value = 0.005313231574131686
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 22.906497
$ws.Cells.Item(2, 14).Value = 68.719491
$ws.Cells.Item(2, 15).Value = 0.9446038650914245
$ws.Cells.Item(2, 16).Value = 0.9446038650914245
$ws.Cells.Item(2, 17).Value = 0.8993701562120001
$ws.Cells.Item(2, 18).Value = 8.094331405908001
$ws.Cells.Item(2, 19).Value = 0.005018899081050585
$ws.Cells.Item(2, 20).Value = 0.005018899081050584

# Row 3: ECs -> FAPs (Adam23-Itga4)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Adam23"
$ws.Cells.Item(3, 3).Value = "Itga4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2.0
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.03926266666666667
$ws.Cells.Item(3, 8).Value = 0.117788
$ws.Cells.Item(3, 9).Value = 0.005313231574131687
$ws.Cells.Item(3, 10).Value = 0.005313231574131686
$ws.Cells.Item(3, 11).Value = 2.0
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1329193333333333
$ws.Cells.Item(3, 14).Value = 0.3987579999999999
$ws.Cells.Item(3, 15).Value = 0.005481244732096839
$ws.Cells.Item(3, 16).Value = 0.005481244732096839
$ws.Cells.Item(3, 17).Value = 0.005218767478222221
$ws.Cells.Item(3, 18).Value = 0.046968907304
$ws.Cells.Item(3, 19).Value = 0.00002912312257611991
$ws.Cells.Item(3, 20).Value = 0.0000291231225761199

# Row 4: ECs -> sCs (Adam23-Itga4)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Adam23"
$ws.Cells.Item(4, 3).Value = "Itga4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2.0
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.03926266666666667
$ws.Cells.Item(4, 8).Value = 0.117788
$ws.Cells.Item(4, 9).Value = 0.005313231574131687
$ws.Cells.Item(4, 10).Value = 0.005313231574131686
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 1.210428333333333
$ws.Cells.Item(4, 14).Value = 3.631285
$ws.Cells.Item(4, 15).Value = 0.04991489017647865
$ws.Cells.Item(4, 16).Value = 0.04991489017647865
$ws.Cells.Item(4, 17).Value = 0.04752464417555556
$ws.Cells.Item(4, 18).Value = 0.42772179758
$ws.Cells.Item(4, 19).Value = 0.000265209370504982
$ws.Cells.Item(4, 20).Value = 0.0002652093705049819

# Row 5: FAPs -> ECs (Adam23-Itga4)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Adam23"
$ws.Cells.Item(5, 3).Value = "Itga4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 4.402094666666667
$ws.Cells.Item(5, 8).Value = 13.206284
$ws.Cells.Item(5, 9).Value = 0.5957147173375057
$ws.Cells.Item(5, 10).Value = 0.5957147173375056
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 22.906497
$ws.Cells.Item(5, 14).Value = 68.719491
$ws.Cells.Item(5, 15).Value = 0.9446038650914245
$ws.Cells.Item(5, 16).Value = 0.9446038650914245
$ws.Cells.Item(5, 17).Value = 100.836568275716
$ws.Cells.Item(5, 18).Value = 907.529114481444
$ws.Cells.Item(5, 19).Value = 0.5627144244888533
$ws.Cells.Item(5, 20).Value = 0.5627144244888532

# Row 6: FAPs -> FAPs (Adam23-Itga4)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Adam23"
$ws.Cells.Item(6, 3).Value = "Itga4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 4.402094666666667
$ws.Cells.Item(6, 8).Value = 13.206284
$ws.Cells.Item(6, 9).Value = 0.5957147173375057
$ws.Cells.Item(6, 10).Value = 0.5957147173375056
$ws.Cells.Item(6, 11).Value = 2.0
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1329193333333333
$ws.Cells.Item(6, 14).Value = 0.3987579999999999
$ws.Cells.Item(6, 15).Value = 0.005481244732096839
$ws.Cells.Item(6, 16).Value = 0.005481244732096839
$ws.Cells.Item(6, 17).Value = 0.5851234883635554
$ws.Cells.Item(6, 18).Value = 5.266111395272
$ws.Cells.Item(6, 19).Value = 0.003265258156238761
$ws.Cells.Item(6, 20).Value = 0.00326525815623876

# Row 7: FAPs -> sCs (Adam23-Itga4)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Adam23"
$ws.Cells.Item(7, 3).Value = "Itga4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 4.402094666666667
$ws.Cells.Item(7, 8).Value = 13.206284
$ws.Cells.Item(7, 9).Value = 0.5957147173375057
$ws.Cells.Item(7, 10).Value = 0.5957147173375056
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 1.210428333333333
$ws.Cells.Item(7, 14).Value = 3.631285
$ws.Cells.Item(7, 15).Value = 0.04991489017647865
$ws.Cells.Item(7, 16).Value = 0.04991489017647865
$ws.Cells.Item(7, 17).Value = 5.32842011054889
$ws.Cells.Item(7, 18).Value = 47.95578099494
$ws.Cells.Item(7, 19).Value = 0.02973503469241362
$ws.Cells.Item(7, 20).Value = 0.02973503469241361

# Row 8: sCs -> ECs (Adam23-Itga4)
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Adam23"
$ws.Cells.Item(8, 3).Value = "Itga4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 2.948244666666667
$ws.Cells.Item(8, 8).Value = 8.844734
$ws.Cells.Item(8, 9).Value = 0.3989720510883627
$ws.Cells.Item(8, 10).Value = 0.3989720510883626
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 22.906497
$ws.Cells.Item(8, 14).Value = 68.719491
$ws.Cells.Item(8, 15).Value = 0.9446038650914245
$ws.Cells.Item(8, 16).Value = 0.9446038650914245
$ws.Cells.Item(8, 17).Value = 67.53395761226601
$ws.Cells.Item(8, 18).Value = 607.8056185103941
$ws.Cells.Item(8, 19).Value = 0.3768705415215207
$ws.Cells.Item(8, 20).Value = 0.3768705415215206

# Row 9: sCs -> FAPs (Adam23-Itga4)
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Adam23"
$ws.Cells.Item(9, 3).Value = "Itga4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 2.948244666666667
$ws.Cells.Item(9, 8).Value = 8.844734
$ws.Cells.Item(9, 9).Value = 0.3989720510883627
$ws.Cells.Item(9, 10).Value = 0.3989720510883626
$ws.Cells.Item(9, 11).Value = 2.0
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1329193333333333
$ws.Cells.Item(9, 14).Value = 0.3987579999999999
$ws.Cells.Item(9, 15).Value = 0.005481244732096839
$ws.Cells.Item(9, 16).Value = 0.005481244732096839
$ws.Cells.Item(9, 17).Value = 0.3918787155968889
$ws.Cells.Item(9, 18).Value = 3.526908440372
$ws.Cells.Item(9, 19).Value = 0.002186863453281959
$ws.Cells.Item(9, 20).Value = 0.002186863453281959

# Row 10: sCs -> sCs (Adam23-Itga4)
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Adam23"
$ws.Cells.Item(10, 3).Value = "Itga4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 2.948244666666667
$ws.Cells.Item(10, 8).Value = 8.844734
$ws.Cells.Item(10, 9).Value = 0.3989720510883627
$ws.Cells.Item(10, 10).Value = 0.3989720510883626
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 1.210428333333333
$ws.Cells.Item(10, 14).Value = 3.631285
$ws.Cells.Item(10, 15).Value = 0.04991489017647865
$ws.Cells.Item(10, 16).Value = 0.04991489017647865
$ws.Cells.Item(10, 17).Value = 3.568638878132223
$ws.Cells.Item(10, 18).Value = 32.11774990319
$ws.Cells.Item(10, 19).Value = 0.01991464611356005
$ws.Cells.Item(10, 20).Value = 0.01991464611356005

